# Add a new "OBSERVACIONES" / "${ObserChequeo}" row at the end of the
# first (only) table in the document, right after the existing
# "REMISORIO" / "${Remisorio}" row.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Word copies the formatting (rPr/pPr) of the row above when a new row
# is appended, which is exactly the formatting the new row needs here.
$newRow = $t.Rows.Add()

$cell1 = $newRow.Cells.Item(1)
$cell2 = $newRow.Cells.Item(2)

# First cell: a single bold run reading "OBSERVACIONES".
$cell1.Range.Text = "OBSERVACIONES"

# Second cell: "${ObserChequeo}" as plain (non-bold) text first, ...
$cell2.Range.Text = '${ObserChequeo}'

# ... then re-insert the paragraph's XML so the "ObserChequeo" word is
# split into its own run and wrapped in spellcheck proofErr markers,
# matching how Word marks an unrecognised camel-case word - while the
# "${" and "}" stay in their own runs around it.
$p2 = $cell2.Range.Paragraphs.Item(1)
$full2 = $p2.Range

$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr>'
$pPr2 = '<w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0" w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr>'

$body2 = $pPr2 + '<w:r>' + $rPr + '<w:t>${</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>ObserChequeo</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr + '<w:t>}</w:t></w:r>'

$xmlFrag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p>' + $body2 + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$full2.InsertXML($xmlFrag2)

Write-Output "Added OBSERVACIONES row"
